# Fruta / hortaliza, semanal
# Insert two new weekly price records at the top of the data block
# (rows 6 and 7), pushing the existing rows 6-27 down to rows 8-29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 6 (Excel copies the formatting of the
# row above, which is what the target workbook shows: D6/D7 keep the
# date-style s="2").
$ws.Rows.Item(6).Resize(2).Insert()

# ---- New row 6: Brooks / Segunda, Region del Maule ----
$ws.Cells.Item(6,1).Value  = 1
$ws.Cells.Item(6,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(6,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(6,4).Value  = 45260
$ws.Cells.Item(6,5).Value  = 15
$ws.Cells.Item(6,6).Value  = "Fruta"
$ws.Cells.Item(6,7).Value  = 100103
$ws.Cells.Item(6,8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(6,9).Value  = 100103001
$ws.Cells.Item(6,10).Value = "Cereza"
$ws.Cells.Item(6,11).Value = "Brooks"
$ws.Cells.Item(6,12).Value = "Segunda"
$ws.Cells.Item(6,13).Value = 180
$ws.Cells.Item(6,14).Value = 14000
$ws.Cells.Item(6,15).Value = 15000
$ws.Cells.Item(6,16).Value = 14556
$ws.Cells.Item(6,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(6,18).Value = "Región del Maule"
$ws.Cells.Item(6,19).Value = 1456
$ws.Cells.Item(6,20).Value = 10

# ---- New row 7: Early Burlat / Primera, Region de O'Higgins ----
$ws.Cells.Item(7,1).Value  = 1
$ws.Cells.Item(7,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(7,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(7,4).Value  = 45260
$ws.Cells.Item(7,5).Value  = 15
$ws.Cells.Item(7,6).Value  = "Fruta"
$ws.Cells.Item(7,7).Value  = 100103
$ws.Cells.Item(7,8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(7,9).Value  = 100103001
$ws.Cells.Item(7,10).Value = "Cereza"
$ws.Cells.Item(7,11).Value = "Early Burlat"
$ws.Cells.Item(7,12).Value = "Primera"
$ws.Cells.Item(7,13).Value = 200
$ws.Cells.Item(7,14).Value = 23000
$ws.Cells.Item(7,15).Value = 25000
$ws.Cells.Item(7,16).Value = 24000
$ws.Cells.Item(7,17).Value = "$/caja 18 kilos"
$ws.Cells.Item(7,18).Value = "Región de O'Higgins"
$ws.Cells.Item(7,19).Value = 1333
$ws.Cells.Item(7,20).Value = 18
